# Apply updated cryptocurrency price/volume data to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.431.63"
$ws.Range("E2").Value = "  +5.54%  "
$ws.Range("D3").Value = "1.809.48"
$ws.Range("E3").Value = "  +4.25%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'316.58"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.5506"
$ws.Range("E7").Value = "  +10.84%  "
$ws.Range("D8").Value = "'0.3859"
$ws.Range("E8").Value = "  +9.32%  "
$ws.Range("D9").Value = "'0.07588"
$ws.Range("E9").Value = "  +4.91%  "
$ws.Range("D10").Value = "'42.96"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").Value = "'1.136"
$ws.Range("E11").Value = "  +7.71%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("E13").Value = "  +6.04%  "
$ws.Range("D14").Value = "'6.229"
$ws.Range("E14").Value = "  +5.03%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.325"
$ws.Range("E15").Value = "  +7.08%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.800.32"
$ws.Range("E16").Value = "  +3.87%  "
$ws.Range("D17").Value = "'91.19"
$ws.Range("E17").Value = "  +5.74%  "
$ws.Range("E18").Value = "  +3.86%  "
$ws.Range("D19").Value = "'0.06472"
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "'17.26"
$ws.Range("E21").Value = "  +4.05%  "
$ws.Range("D22").Value = "'5.986"
$ws.Range("E22").Value = "  +4.65%  "
$ws.Range("D23").Value = "28.451.59"
$ws.Range("E23").Value = "  +5.33%  "
$ws.Range("D24").Value = "'11.30"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "'2.123"
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("D26").Value = "'157.46"
$ws.Range("E26").Value = "  +2.61%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.447"
$ws.Range("E27").Value = "  +15.95%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.70"
$ws.Range("E28").Value = "  +4.28%  "
$ws.Range("D29").Value = "2.020.01"
$ws.Range("E29").Value = "  +4.27%  "
$ws.Range("D30").Value = "'123.83"
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").Value = "'1.170"
$ws.Range("E31").Value = "  +10.66%  "
$ws.Range("D32").Value = "'0.1035"
$ws.Range("E32").Value = "  +9.38%  "
$ws.Range("D33").Value = "'5.761"
$ws.Range("E33").Value = "  +7.34%  "
$ws.Range("D34").Value = "'3.649"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("D35").Value = "'0.2277"
$ws.Range("E35").Value = "  +14.38%  "
$ws.Range("D36").Value = "'8.949"
$ws.Range("E36").Value = "  +19.88%  "
$ws.Range("D37").Value = "'0.02333"
$ws.Range("E37").Value = "  +6.61%  "
$ws.Range("D38").Value = "'0.06251"
$ws.Range("E38").Value = "  +5.53%  "
$ws.Range("D39").Value = "'11.61"
$ws.Range("E39").Value = "  +5.62%  "
$ws.Range("D40").Value = "'0.6390"
$ws.Range("E40").Value = "  +6.76%  "
$ws.Range("D41").Value = "'5.021"
$ws.Range("E41").Value = "  +5.88%  "
$ws.Range("D42").Value = "'1.180"
$ws.Range("E42").Value = "  +6.65%  "
$ws.Range("D43").Value = "'0.9999"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'1.386"
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("D45").Value = "'13.43"
$ws.Range("E45").Value = "  +4.63%  "
$ws.Range("D46").Value = "'0.6007"
$ws.Range("E46").Value = "  +6.82%  "
$ws.Range("D47").Value = "'3.689"
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("D48").Value = "'123.26"
$ws.Range("E48").Value = "  +3.41%  "
$ws.Range("D49").Value = "'1.975"
$ws.Range("E49").Value = "  +6.91%  "
$ws.Range("D50").Value = "'1.145"
$ws.Range("E50").Value = "  +4.40%  "
$ws.Range("D51").Value = "'0.06929"
$ws.Range("E51").Value = "  +4.07%  "
